$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) sometimes contains values that look like plain
# numbers (e.g. "578.00", "0.0000184"). Excel auto-detects these and
# would silently convert them to numeric cells, stripping formatting
# such as trailing zeros. Forcing a Text number format right before
# writing those specific cells keeps them as text, matching the source
# data which stores every Price/Volume cell as a text string.
$ws.Range("D2").Value = '63.472.31'
$ws.Range("E2").Value = '  +5.85%  '
$ws.Range("D3").Value = '3.399.48'
$ws.Range("E3").Value = '  +6.71%  '
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.00'
$ws.Range("E5").Value = '  +7.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.62'
$ws.Range("E6").Value = '  +7.33%  '
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").Value = '3.405.76'
$ws.Range("E8").Value = '  +6.76%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.534'
$ws.Range("E9").Value = '  +0.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.48'
$ws.Range("E10").Value = '  +1.89%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.121'
$ws.Range("E11").Value = '  +7.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.436'
$ws.Range("E12").Value = '  +0.91%  '
$ws.Range("D13").Value = '3.979.11'
$ws.Range("E13").Value = '  +6.48%  '
$ws.Range("E14").Value = '  +0.43%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000184'
$ws.Range("E15").Value = '  +7.05%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '27.11'
$ws.Range("E16").Value = '  +5.20%  '
$ws.Range("D17").Value = '63.549.82'
$ws.Range("E17").Value = '  +5.96%  '
$ws.Range("D18").Value = '3.384.06'
$ws.Range("E18").Value = '  +5.98%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.40'
$ws.Range("E19").Value = '  +2.64%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.02'
$ws.Range("E20").Value = '  +5.92%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.44'
$ws.Range("E21").Value = '  +3.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '388.16'
$ws.Range("E22").Value = '  +5.27%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.999'
$ws.Range("E23").Value = '  -0.03%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.536'
$ws.Range("E24").Value = '  +2.48%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '70.69'
$ws.Range("E25").Value = '  +1.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.58'
$ws.Range("E26").Value = '  +12.37%  '
$ws.Range("B27").Value = 'PEPE'
$ws.Range("C27").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000104'
$ws.Range("E27").Value = '  +19.00%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.180'
$ws.Range("E28").Value = '  +6.59%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.999'
$ws.Range("E29").Value = '  -0.22%  '
$ws.Range("E30").Value = '  +7.88%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.48'
$ws.Range("E31").Value = '  +6.13%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.34'
$ws.Range("E32").Value = '  +12.69%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '23.14'
$ws.Range("E33").Value = '  +2.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.63'
$ws.Range("E34").Value = '  +6.75%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.72'
$ws.Range("E35").Value = '  +2.12%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.50'
$ws.Range("E36").Value = '  +10.28%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '158.22'
$ws.Range("E37").Value = '  +0.66%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.88'
$ws.Range("E38").Value = '  +11.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.60'
$ws.Range("E39").Value = '  +4.86%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0760'
$ws.Range("E40").Value = '  +7.74%  '
$ws.Range("D41").Value = '2.905.81'
$ws.Range("E41").Value = '  +4.47%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0323'
$ws.Range("E42").Value = '  +3.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.765'
$ws.Range("E43").Value = '  +6.44%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.41'
$ws.Range("E44").Value = '  +4.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.31'
$ws.Range("E45").Value = '  +2.09%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.07'
$ws.Range("E46").Value = '  +8.90%  '
$ws.Range("D47").Value = '3.445.01'
$ws.Range("E47").Value = '  +6.73%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.29'
$ws.Range("E48").Value = '  +8.22%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '300.26'
$ws.Range("E49").Value = '  +14.37%  '
$ws.Range("E50").Value = '  -0.89%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.33'
$ws.Range("E51").Value = '  +3.07%  '
